$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the specified rows per repulled data
$ws.Range("F2").Value = -9
$ws.Range("F3").Value = 4
$ws.Range("F6").Value = 3
$ws.Range("F8").Value = -6
$ws.Range("F9").Value = -1
$ws.Range("F15").Value = 4
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = -1
$ws.Range("F24").Value = 1
$ws.Range("F25").Value = -3
$ws.Range("F26").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("F35").Value = -3
$ws.Range("F41").Value = 3
$ws.Range("F42").Value = 1
$ws.Range("F43").Value = 1
$ws.Range("F47").Value = 0
$ws.Range("F51").Value = -5
